# C4-PowerPoint.pptx — swap the theme colour scheme used by the slide
# master from "Integral" to the default "Office Theme" palette
# (ppt/theme/theme1.xml: a:clrScheme Integral -> Office).
#
# Helper: turn an RRGGBB hex string into the decimal value PowerPoint's
# ColorFormat.RGB / VBA RGB() uses (0x00BBGGRR, i.e. R + G*256 + B*65536).
function HexToRgbVal([string]$hex) {
    $v = [Convert]::ToInt32($hex, 16)
    $r = ($v -shr 16) -band 0xFF
    $g = ($v -shr 8) -band 0xFF
    $b = $v -band 0xFF
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

# Target palette: the built-in Office Theme colour scheme.
# Index order matches ThemeColorScheme.Item(1..12):
#   1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#   5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink
$officeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Item($i).RGB = HexToRgbVal($officeColors[$i - 1])
}
